# The deck ships two theme parts: the slide master's theme (the "Integral"
# palette) and a second theme used only by the notes master (the stock
# "Office Theme" palette). The edit swaps which palette lives in which part:
# the slide master's theme becomes the Office palette, and the notes-master
# theme becomes the Integral palette.
#
# PowerPoint's object model only ever surfaces a single, editable theme -
# the one attached to ActivePresentation.SlideMaster - so we recolor that
# theme's 12-color scheme to the target ("Office Theme") values one swatch
# at a time via ThemeColorScheme. Note: COM's ColorFormat.RGB uses the
# 0x00BBGGRR long form, not 0x00RRGGBB.

function Set-ThemeRGB($colorScheme, $index, $r, $g, $b) {
    $colorScheme.Item($index).RGB = ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$cs = $theme.ThemeColorScheme

# dk1 / lt1 stay black / white in both palettes - set for completeness.
Set-ThemeRGB $cs 1  0x00 0x00 0x00   # dk1
Set-ThemeRGB $cs 2  0xFF 0xFF 0xFF   # lt1

# Remaining swatches move from the "Integral" palette to the "Office Theme" palette.
Set-ThemeRGB $cs 3  0x44 0x54 0x6A   # dk2
Set-ThemeRGB $cs 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeRGB $cs 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeRGB $cs 6  0xED 0x7D 0x31   # accent2
Set-ThemeRGB $cs 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeRGB $cs 8  0xFF 0xC0 0x00   # accent4
Set-ThemeRGB $cs 9  0x44 0x72 0xC4   # accent5
Set-ThemeRGB $cs 10 0x70 0xAD 0x47   # accent6
Set-ThemeRGB $cs 11 0x05 0x63 0xC1   # hlink
Set-ThemeRGB $cs 12 0x95 0x4F 0x72   # folHlink
